$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fix the clinician name recorded for the first two appointments: the
# placeholder "Provider Test" is replaced with the real clinician name.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "John Doe"
$ws.Range("E3").Value = "John Doe"

# ---------------------------------------------------------------------------
# New data rows (4-10): appointments for Patient3..Patient9, alternating
# clinicians "John Doe" / "Bob Doe".
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=4;  Date=45294; Code=90833; Patient="Patient3"; Clinician="Bob Doe"  },
    @{ Row=5;  Date=45295; Code=90834; Patient="Patient4"; Clinician="John Doe" },
    @{ Row=6;  Date=45296; Code=90835; Patient="Patient5"; Clinician="Bob Doe"  },
    @{ Row=7;  Date=45297; Code=90836; Patient="Patient6"; Clinician="Bob Doe"  },
    @{ Row=8;  Date=45298; Code=90837; Patient="Patient7"; Clinician="John Doe" },
    @{ Row=9;  Date=45299; Code=90838; Patient="Patient8"; Clinician="John Doe" },
    @{ Row=10; Date=45300; Code=90839; Patient="Patient9"; Clinician="John Doe" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.Code
    $ws.Range("C$row").Value = "Test"
    $ws.Range("D$row").Value = $r.Patient
    $ws.Range("E$row").Value = $r.Clinician
    $ws.Range("F$row").Value = 75
    $ws.Range("G$row").Value = 50
    $ws.Range("H$row").Value = 150
    $ws.Range("I$row").Value = 100
    $ws.Range("J$row").Value = "Appointment"
}

# ---------------------------------------------------------------------------
# Apply the same Date / Service-Code number formatting used in rows 2-3 to
# the newly added rows, as well as to a handful of extra blank rows
# (11-15) that only carry that formatting and no data yet.
# ---------------------------------------------------------------------------
$ws.Range("A2:B3").Copy()
$ws.Range("A4:B15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Widen column E (Clinician Name) now that it holds longer values.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 17.7109375
